# 141: 31/12 09:34 (LP1912 + 6203+6173) scrape run appended.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912": new scrape timestamp, updated row count, 8 new arrival rows.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 31/12/2025 06:34:40"
$ws1.Range("A3").Value = "Total filas: 692"

$newRows1 = @(
    @{ r = 686; B = "06:34:29"; C = "06:54"; D = "14_ABASTO";     E = 20; F = "LP1912"; G = "31/12/2025" },
    @{ r = 687; B = "06:34:29"; C = "07:01"; D = "16_SANTA ANA";  E = 27; F = "LP1912"; G = "31/12/2025" },
    @{ r = 688; B = "06:34:29"; C = "07:13"; D = "14X44_ABASTO";  E = 39; F = "LP1912"; G = "31/12/2025" },
    @{ r = 689; B = "06:34:29"; C = "07:16"; D = "16_SANTA ANA";  E = 42; F = "LP1912"; G = "31/12/2025" },
    @{ r = 690; B = "06:34:29"; C = "07:51"; D = "15_ABASTO";     E = 77; F = "LP1912"; G = "31/12/2025" },
    @{ r = 691; B = "06:34:29"; C = "08:02"; D = "23_HERNANDEZ";  E = 88; F = "LP1912"; G = "31/12/2025" },
    @{ r = 692; B = "06:34:29"; C = "08:03"; D = "17_ROMERO";     E = 89; F = "LP1912"; G = "31/12/2025" },
    @{ r = 693; B = "06:34:29"; C = "08:13"; D = "10_OLMOS";      E = 99; F = "LP1912"; G = "31/12/2025" }
)

foreach ($row in $newRows1) {
    $ws1.Cells.Item($row.r, 2).Value = $row.B
    $ws1.Cells.Item($row.r, 3).Value = $row.C
    $ws1.Cells.Item($row.r, 4).Value = $row.D
    $ws1.Cells.Item($row.r, 5).Value = $row.E
    $ws1.Cells.Item($row.r, 6).Value = $row.F
    $ws1.Cells.Item($row.r, 7).Value = $row.G
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215": only the scrape timestamp moves; no new rows here.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 31/12/2025 06:34:40"

# ---------------------------------------------------------------------------
# Sheet "6203-6173": new scrape timestamp, updated row count, 1 new arrival row.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 31/12/2025 06:34:40"
$ws3.Range("A3").Value = "Total filas: 83"

$ws3.Cells.Item(84, 2).Value = "31/12/2025"
$ws3.Cells.Item(84, 3).Value = "06:34:39"
$ws3.Cells.Item(84, 4).Value = "07:27"
$ws3.Cells.Item(84, 5).Value = "215A_LA PLATA"
$ws3.Cells.Item(84, 6).Value = 53
$ws3.Cells.Item(84, 7).Value = "L6173"
